$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 82, shifting existing rows 82:97 down to 83:98
$ws.Rows.Item(82).Insert()

# Populate the newly inserted row 82 with the new weekly price entry
$ws.Range("A82").Value = 4
$ws.Range("B82").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C82").Value = "Los Lagos"
$ws.Range("D82").Value = 44551
$ws.Range("E82").Value = 10
$ws.Range("F82").Value = "Fruta"
$ws.Range("G82").Value = 100103
$ws.Range("H82").Value = "Frutos de hueso (carozo)"
$ws.Range("I82").Value = 100103002
$ws.Range("J82").Value = "Ciruela"
$ws.Range("K82").Value = "Red Beaut"
$ws.Range("L82").Value = "Primera"
$ws.Range("M82").Value = 700
$ws.Range("N82").Value = 19000
$ws.Range("O82").Value = 20000
$ws.Range("P82").Value = 19500
$ws.Range("Q82").Value = "`$/caja 15 kilos granel"
$ws.Range("R82").Value = "Región Metropolitana"
$ws.Range("S82").Value = 1300
$ws.Range("T82").Value = 15
